$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the Treatment tab query text in cell B5 ---
# Remove the redundant CONCAT() wrapper around REPLACE(...)
$old = $ws.Range("B5").Value2
$new = $old.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$ws.Range("B5").Value2 = $new

# --- Update the sheet view ---
# Scroll back so the top-left visible cell is A1 (removes topLeftCell="A4")
# and change the active selection from C5 to B2.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
